# PROS-13075 - CCRU - POS KPI 2020 change
#
# The KPI rows (2-14) get re-sorted alphabetically by KPI name (columns
# B, D, E and F all carry the same KPI-name text per row), and the
# "and-or" typo in the cooler-doors KPI is fixed to "and/or" at the same
# time. Row 14's height goes back to the default 15pt, the A/B/C/D:F
# column widths are tightened, and the selection/active cell move to
# A1:F14 / B1 (as after a manual sort + reselect of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (alphabetically sorted) KPI names for rows 2..14, fixing the
# "and-or" -> "and/or" typo along the way.
$values = @(
    "CCH coolers quality",
    "CCH coolers quality (Prime Pos/Max15/Merch STD/Occupancy/Lights&chilled)",
    "CCH products present in Customers menu",
    "CCH shelf share in Energy",
    "CCH shelf share in Juice",
    "CCH shelf share in SSD",
    "CCH shelf share in Tea",
    "CCH shelf share in Water",
    "Number of CCH activation points in NARTD",
    "Number of CCH cooler doors and/or equivalent in Customer coolers",
    "Number of CCH displays points of interaction",
    "Number of NCB core assortment available in-store",
    "Number of SSD core assortment available in-store"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $text = $values[$i]
    $ws.Range("B$row").Value = $text
    $ws.Range("D$row").Value = $text
    $ws.Range("E$row").Value = $text
    $ws.Range("F$row").Value = $text
}

# Row 14 loses its custom (13.8pt) height, going back to the sheet
# default of 15pt.
$ws.Range("A14:F14").RowHeight = 15

# Column widths were tightened (B/D:F down to a round 61 characters).
# The ColumnWidth setter here snaps to a 1/6-character pixel grid, so
# feed it the pre-image of that rounding for each target width.
$ws.Range("A:A").ColumnWidth = 15.8333333333333
$ws.Range("B:B").ColumnWidth = 60.1666666666667
$ws.Range("C:C").ColumnWidth = 17
$ws.Range("D:F").ColumnWidth = 60.1666666666667

# Reselect the table: whole A1:F14 block, active cell on B1.
$ws.Range("A1:F14").Select() | Out-Null
